# English/Spanish translation workbook: add new localisation rows
# (deleting-report spinner/errors, hours/payroll/pay period headers,
# "Messages from Tino" inbox strings) and widen the two text columns
# so the longer English/Spanish strings are readable.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths: columns B and C no longer share one width ---
$ws.Columns("B").ColumnWidth = 123.42578125
$ws.Columns("C").ColumnWidth = 115.7109375

# --- Row height bookkeeping noticed while editing ---
$ws.Rows(7).RowHeight = 15
$ws.Rows(8).RowHeight = 409.5

# --- New translation rows appended after the existing A2:C85 table ---
$newRows = @(
    @("spinner_deleting_report", "Deleting report...", "Eliminar informe..."),
    @("error_deleting_report_message", "Error deleting report. Please try again later.", "Error al eliminar el informe. Por favor, inténtelo de nuevo más tarde."),
    @("error_server_connect_message", "Could not connect to server. Please try again later.", "No se pudo conectar al servidor. Por favor, inténtelo de nuevo más tarde."),
    @("error_fetching_reports_title", "Connection Error", "Error de Conexión"),
    @("error_fetching_reports_message", "Could not connect to server to retrieve reports. Please try again later.", "No se pudo conectar al servidor para recuperar informes. Por favor, inténtelo de nuevo más tarde."),
    @("hours_header", "Hours", "Horas"),
    @("payroll_period", "Payroll Period", "Período de Nómina"),
    @("pay_period", "Pay Period", "Período de Pago"),
    @("messages_title", "Messages from Tino", "Mensajes de Tino"),
    @("message_date", "Date", "Fecha"),
    @("message_from", "Sender", "Remitente"),
    @("message_subject", "Subject", "Tema"),
    @("message_done", "Done", "Hecho"),
    @("messages_no_messages", "No messages available", "No hay mensajes disponibles")
)

$startRow = 86
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $newRows[$i]

    $cellA = $ws.Range("A$r")
    $cellBC = $ws.Range("B$r`:C$r")

    $cellA.Interior.Color = 65535
    $cellBC.Interior.Color = 65535
    $cellBC.WrapText = $true

    $ws.Range("A$r").Value = $rowVals[0]
    $ws.Range("B$r").Value = $rowVals[1]
    $ws.Range("C$r").Value = $rowVals[2]
}

# --- Selection / scroll position left after the edits ---
$ws.Range("A86:C99").Select()

Write-Output "applied translation updates"
